$d = $word.ActiveDocument

# 1. Merge "L" + "O-umlaut" + "VE" runs into a single "LÖVE" run.
$d.Content.Find.Execute("LÖVE", $false, $false, $false, $false, $false, `
    $true, 1, $false, "LÖVE", 2) | Out-Null

# 2. Remove the stray _GoBack bookmark that currently sits after "OpenMP".
$d.Bookmarks("_GoBack").Delete()

# 3. Trim the Physics Senior Thesis bullet.
$d.Content.Find.Execute( `
    "collaboration with a professor, building on an existing project, to increase efficiency by an order of magnitude.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "increase efficiency of existing algorithm by an order of magnitude", 2) | Out-Null

# 4. Reword the Huckster bullet.
$d.Content.Find.Execute( `
    "Huckster – fan-made hero class for the boardgame Shadows of Brimstone", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Huckster – fan-made hero for the board game Shadows of Brimstone", 2) | Out-Null

# 5. Merge the hyperlink display-text runs back into one run.
$d.Content.Find.Execute( `
    "https://boardgamegeek.com/filepage/116086/new-hero-huckster", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "https://boardgamegeek.com/filepage/116086/new-hero-huckster", 2) | Out-Null

# 6. Add "University" into the Oarsman line.
$d.Content.Find.Execute( `
    "Oarsman for San Diego State Men’s Crew (2012-2014)", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Oarsman for San Diego State University Men’s Crew (2012-2014)", 2) | Out-Null

# 7. Remove the "5:00 am practice..." bullet paragraph entirely -- this
#    leaves the following (already-empty) paragraph as the new home for
#    the document's _GoBack bookmark, matching real Word's behaviour of
#    re-anchoring _GoBack to the site of the most recent edit.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "5:00 am practice, 5 days a week, 30+ hours/week training`r") {
        $targetPara = $p
    }
}
$deleteStart = $targetPara.Range.Start
$targetPara.Range.Delete()

# 8. Re-create _GoBack collapsed at the start of the paragraph that used
#    to follow the deleted bullet (now sitting right where the bullet
#    used to start).
$bmPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -eq $deleteStart) {
        $bmPara = $p
    }
}
$bmRange = $bmPara.Range.Duplicate()
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange)
